# Applies numeric corrections scraped from the scheduled-runner commit.
$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3509
$ws.Range("J32").Value = 3408.111
$ws.Range("L32").Value = 3408.111
$ws.Range("N32").Value = -4060.111
$ws.Range("H33").Value = 110.3
$ws.Range("I33").Value = 105.888885
$ws.Range("K33").Value = 105.888885
$ws.Range("M33").Value = 123.111115
$ws.Range("H113").Value = 7159.2144
$ws.Range("I113").Value = 5525
$ws.Range("J113").Value = 7812.9
$ws.Range("K113").Value = 5525
$ws.Range("L113").Value = 7812.9
$ws.Range("M113").Value = -2271
$ws.Range("N113").Value = -14320.9
$ws.Range("H138").Value = 5456.483
$ws.Range("I138").Value = 4379.636
$ws.Range("J138").Value = 6114.5557
$ws.Range("K138").Value = 13138.908
$ws.Range("L138").Value = 18343.6671
$ws.Range("M138").Value = -7998.908000000001
$ws.Range("N138").Value = -28623.6671

# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2596.9778
$ws.Range("I2").Value = 1873.1666
$ws.Range("K2").Value = 1873.1666
$ws.Range("M2").Value = -1760.1666
$ws.Range("H32").Value = 1763981.8
$ws.Range("I32").Value = 1859513.8
$ws.Range("K32").Value = 1859513.8
$ws.Range("M32").Value = -1859226.8
$ws.Range("H45").Value = 3442.1428
$ws.Range("I45").Value = 4155.6665
$ws.Range("J45").Value = 2157.8
$ws.Range("K45").Value = 4155.6665
$ws.Range("L45").Value = 2157.8
$ws.Range("M45").Value = -3778.6665
$ws.Range("N45").Value = -2911.8
$ws.Range("H61").Value = 9652.200000000001
$ws.Range("I61").Value = 5478.8
$ws.Range("K61").Value = 5478.8
$ws.Range("M61").Value = -5266.8
$ws.Range("H74").Value = 4633631.5
$ws.Range("I74").Value = 7814711
$ws.Range("K74").Value = 7814711
$ws.Range("M74").Value = -7813837
$ws.Range("H77").Value = 4633631.5
$ws.Range("I77").Value = 7814711
$ws.Range("K77").Value = 39073555
$ws.Range("M77").Value = -39069187
$ws.Range("H116").Value = 2596.9778
$ws.Range("I116").Value = 1873.1666
$ws.Range("K116").Value = 1873.1666
$ws.Range("M116").Value = 420.8334
$ws.Range("H136").Value = 9652.200000000001
$ws.Range("I136").Value = 5478.8
$ws.Range("K136").Value = 16436.4
$ws.Range("M136").Value = -13886.4

# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2596.9778
$ws.Range("I3").Value = 1873.1666
$ws.Range("K3").Value = 1873.1666
$ws.Range("M3").Value = -1759.1666
$ws.Range("H134").Value = 613768.5
$ws.Range("I134").Value = 746061.8
$ws.Range("K134").Value = 2238185.4
$ws.Range("M134").Value = -2235650.4
$ws.Range("H139").Value = 80656.8
$ws.Range("J139").Value = 80656.8
$ws.Range("L139").Value = 80656.8
$ws.Range("N139").Value = -90936.8

# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 11441.632
$ws.Range("I31").Value = 4058
$ws.Range("J31").Value = 15748.75
$ws.Range("K31").Value = 4058
$ws.Range("L31").Value = 15748.75
$ws.Range("M31").Value = -3763
$ws.Range("N31").Value = -16338.75
$ws.Range("H34").Value = 11441.632
$ws.Range("I34").Value = 4058
$ws.Range("J34").Value = 15748.75
$ws.Range("K34").Value = 4058
$ws.Range("L34").Value = 15748.75
$ws.Range("M34").Value = -3856
$ws.Range("N34").Value = -16152.75
$ws.Range("H58").Value = 776792.9
$ws.Range("I58").Value = 1034885.3
$ws.Range("J58").Value = 2515.5
$ws.Range("K58").Value = 1034885.3
$ws.Range("L58").Value = 2515.5
$ws.Range("M58").Value = -1034682.3
$ws.Range("N58").Value = -2921.5
$ws.Range("H132").Value = 19263470
$ws.Range("I132").Value = 51391.285
$ws.Range("K132").Value = 154173.855
$ws.Range("M132").Value = -151643.855
$ws.Range("H136").Value = 776792.9
$ws.Range("I136").Value = 1034885.3
$ws.Range("J136").Value = 2515.5
$ws.Range("K136").Value = 3104655.9
$ws.Range("L136").Value = 7546.5
$ws.Range("M136").Value = -3102105.9
$ws.Range("N136").Value = -12646.5

# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H74").Value = 6983.25
$ws.Range("I74").Value = 3997.5
$ws.Range("K74").Value = 11992.5
$ws.Range("M74").Value = -10931.5
$ws.Range("H77").Value = 6983.25
$ws.Range("I77").Value = 3997.5
$ws.Range("K77").Value = 35977.5
$ws.Range("M77").Value = -30673.5
$ws.Range("H87").Value = 18927.334
$ws.Range("I87").Value = 6750
$ws.Range("J87").Value = 25016
$ws.Range("K87").Value = 20250
$ws.Range("L87").Value = 75048
$ws.Range("M87").Value = -19002
$ws.Range("N87").Value = -77544
$ws.Range("H90").Value = 18927.334
$ws.Range("I90").Value = 6750
$ws.Range("J90").Value = 25016
$ws.Range("K90").Value = 60750
$ws.Range("L90").Value = 225144
$ws.Range("M90").Value = -54510
$ws.Range("N90").Value = -237624
$ws.Range("H114").Value = 2315.3845
$ws.Range("I114").Value = 776.93335
$ws.Range("K114").Value = 2330.80005
$ws.Range("M114").Value = 923.1999500000002
$ws.Range("H131").Value = 13824
$ws.Range("J131").Value = 17255.133
$ws.Range("L131").Value = 51765.399
$ws.Range("N131").Value = -61845.399
$ws.Range("H140").Value = 2306.8235
$ws.Range("I140").Value = 2306.8235
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 6920.470499999999
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -1740.470499999999
$ws.Range("H141").Value = 2906.25
$ws.Range("I141").Value = 2906.25
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 8718.75
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -3538.75
$ws.Range("N140").ClearContents()
$ws.Range("N141").ClearContents()

# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 458172
$ws.Range("I80").Value = 502989.2
$ws.Range("J80").Value = 10000
$ws.Range("K80").Value = 502989.2
$ws.Range("L80").Value = 10000
$ws.Range("M80").Value = -501991.2
$ws.Range("N80").Value = -11996
$ws.Range("H83").Value = 458172
$ws.Range("I83").Value = 502989.2
$ws.Range("J83").Value = 10000
$ws.Range("K83").Value = 2514946
$ws.Range("L83").Value = 50000
$ws.Range("M83").Value = -2509954
$ws.Range("N83").Value = -59984
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("H102").Value = 2541.853
$ws.Range("I102").Value = 2332.96
$ws.Range("K102").Value = 2332.96
$ws.Range("M102").Value = -710.96
$ws.Range("H132").Value = 4919.7
$ws.Range("I132").Value = 4931.6665
$ws.Range("J132").Value = 4901.75
$ws.Range("K132").Value = 14794.9995
$ws.Range("L132").Value = 14705.25
$ws.Range("M132").Value = -12264.9995
$ws.Range("N132").Value = -19765.25
$ws.Range("N87").ClearContents()
$ws.Range("N90").ClearContents()

# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7268.077
$ws.Range("I7").Value = 6748.9
$ws.Range("J7").Value = 8998.666999999999
$ws.Range("K7").Value = 6748.9
$ws.Range("L7").Value = 8998.666999999999
$ws.Range("M7").Value = -6636.9
$ws.Range("N7").Value = -9222.666999999999
$ws.Range("H82").Value = 1188.4
$ws.Range("I82").Value = 826
$ws.Range("J82").Value = 2034
$ws.Range("K82").Value = 826
$ws.Range("L82").Value = 2034
$ws.Range("M82").Value = -465
$ws.Range("N82").Value = -2756
$ws.Range("H85").Value = 1188.4
$ws.Range("I85").Value = 826
$ws.Range("J85").Value = 2034
$ws.Range("K85").Value = 826
$ws.Range("L85").Value = 2034
$ws.Range("M85").Value = 422
$ws.Range("N85").Value = -4530
$ws.Range("H93").Value = 1352.4286
$ws.Range("I93").Value = 1120.05
$ws.Range("J93").Value = 6000
$ws.Range("K93").Value = 1120.05
$ws.Range("L93").Value = 6000
$ws.Range("M93").Value = 127.95
$ws.Range("N93").Value = -8496
$ws.Range("H122").Value = 73809.8
$ws.Range("J122").Value = 211939.8
$ws.Range("L122").Value = 635819.3999999999
$ws.Range("N122").Value = -640719.3999999999
$ws.Range("H126").Value = 7268.077
$ws.Range("I126").Value = 6748.9
$ws.Range("J126").Value = 8998.666999999999
$ws.Range("K126").Value = 20246.7
$ws.Range("L126").Value = 26996.001
$ws.Range("M126").Value = -17776.7
$ws.Range("N126").Value = -31936.001
$ws.Range("H132").Value = 1445502.8
$ws.Range("I132").Value = 1924731.1
$ws.Range("K132").Value = 5774193.300000001
$ws.Range("M132").Value = -5771663.300000001
$ws.Range("H136").Value = 4584.909
$ws.Range("J136").Value = 7149.6665
$ws.Range("L136").Value = 21448.9995
$ws.Range("N136").Value = -26548.9995

# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 2527.5
$ws.Range("I113").Value = 2440.3333
$ws.Range("J113").Value = 2579.8
$ws.Range("K113").Value = 7320.999899999999
$ws.Range("L113").Value = 7739.400000000001
$ws.Range("M113").Value = -5150.999899999999
$ws.Range("N113").Value = -12079.4
$ws.Range("H122").Value = 3144.743
$ws.Range("I122").Value = 2845.8438
$ws.Range("K122").Value = 8537.5314
$ws.Range("M122").Value = -6087.5314
$ws.Range("H132").Value = 142283250
$ws.Range("I132").Value = 34563456
$ws.Range("K132").Value = 103690368
$ws.Range("M132").Value = -103687838
